$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A58").Value = "1b650324aafac0e6e6f0c473eccff258"
$ws.Range("B58").Value = "Homo sapiens"
$ws.Range("C58").Value = "Human"
$ws.Range("D58").Value = "Human"
$ws.Range("A59").Value = "407f08a29007a8a153222d82ef47d408"
$ws.Range("B59").Value = "Menidia menidia"
$ws.Range("C59").Value = "Atlantic silverside"
$ws.Range("D59").Value = "Teleost Fish"
$ws.Range("A76").Value = "09351b480c58a99c4150d54ebbc97c6c"
$ws.Range("B76").Value = "Menidia menidia"
$ws.Range("C76").Value = "Atlantic silverside"
$ws.Range("D76").Value = "Teleost Fish"
$ws.Range("A77").Value = "d34820a8c9954e292ea9dbc76f4275b4"
$ws.Range("B77").Value = "Homo sapiens"
$ws.Range("C77").Value = "Human"
$ws.Range("D77").Value = "Human"
$ws.Range("A93").Value = "680475954df3011ebba1033f1b2f2a86"
$ws.Range("B93").Value = "Prionotus carolinus"
$ws.Range("C93").Value = "Northern sea robin"
$ws.Range("D93").Value = "Teleost Fish"
$ws.Range("A94").Value = "4c5905c5ab539613d9c3069d0ae54188"
$ws.Range("B94").Value = "Bos taurus"
$ws.Range("C94").Value = "Cow"
$ws.Range("D94").Value = "Livestock"
$ws.Range("A95").Value = "db8615250f29272019fe417d96bf08f3"
$ws.Range("B95").Value = "Brevoortia tyrannus"
$ws.Range("C95").Value = "Atlantic menhaden"
$ws.Range("D95").Value = "Teleost Fish"
$ws.Range("A105").Value = "f5e0ea6fe3e45da9605b758c440ae692"
$ws.Range("B105").Value = "Brevoortia tyrannus"
$ws.Range("C105").Value = "Atlantic menhaden"
$ws.Range("D105").Value = "Teleost Fish"
$ws.Range("A106").Value = "6d1668646cf923fa90217b0797de7a7d"
$ws.Range("B106").Value = "Pseudopleuronectes americanus or Myzopsetta ferruginea"
$ws.Range("C106").Value = "Winter or Yellowtail flounder"
$ws.Range("D106").Value = "Teleost Fish"
$ws.Range("A113").Value = "5e733a21f67e541f28ed4bf4fe025044"
$ws.Range("B113").Value = "Paralichthys dentatus"
$ws.Range("C113").Value = "Summer flounder"
$ws.Range("D113").Value = "Teleost Fish"
$ws.Range("A114").Value = "f524c4b860dec1e6b994c28dd8e4b75e"
$ws.Range("B114").Value = "Ammodytes dubius"
$ws.Range("C114").Value = "Northern sand lance"
$ws.Range("D114").Value = "Teleost Fish"
$ws.Range("A144").Value = "9db3dc01519672b43908456a37b27b4d"
$ws.Range("B144").Value = "Fundulus heteroclitus"
$ws.Range("C144").Value = "Mummichog"
$ws.Range("D144").Value = "Teleost Fish"
$ws.Range("A145").Value = "0df37a1b74088f0e7410a1b78cada881"
$ws.Range("B145").Value = "Engraulis eurystole"
$ws.Range("C145").Value = "Silver anchovy"
$ws.Range("D145").Value = "Teleost Fish"
$ws.Range("A150").Value = "0f011be680aec3ee4b12b1b139902251"
$ws.Range("B150").Value = "Apeltes quadracus"
$ws.Range("C150").Value = "Fourspine stickleback"
$ws.Range("D150").Value = "Teleost Fish"
$ws.Range("A151").Value = "94a944154183c458facbab20fe39ffa9"
$ws.Range("B151").Value = "Apeltes quadracus"
$ws.Range("C151").Value = "Fourspine stickleback"
$ws.Range("D151").Value = "Teleost Fish"
$ws.Range("A172").Value = "191ed810bb884ed43fa1919f6da3d82a"
$ws.Range("B172").Value = "Pseudopleuronectes americanus or Myzopsetta ferruginea"
$ws.Range("C172").Value = "Winter or Yellowtail flounder"
$ws.Range("D172").Value = "Teleost Fish"
$ws.Range("A173").Value = "148aa3594130e12c353383f68bfa0b6a"
$ws.Range("B173").Value = "Apeltes quadracus"
$ws.Range("C173").Value = "Fourspine stickleback"
$ws.Range("D173").Value = "Teleost Fish"
$ws.Range("A174").Value = "c73cefb2b4ac8de08ae0c68341cbb28f"
$ws.Range("B174").Value = "Ammodytes dubius"
$ws.Range("C174").Value = "Northern sand lance"
$ws.Range("D174").Value = "Teleost Fish"
$ws.Range("A175").Value = "ed00c11476e9a07d3441cb0a1073d3ab"
$ws.Range("B175").Value = "Etropus microstomus"
$ws.Range("C175").Value = "Smallmouth flounder"
$ws.Range("D175").Value = "Teleost Fish"
$ws.Range("A181").Value = "ff405ebc8992c59ba51a99e33a12fe74"
$ws.Range("B181").Value = "Pseudopleuronectes americanus or Myzopsetta ferruginea"
$ws.Range("C181").Value = "Winter or Yellowtail flounder"
$ws.Range("D181").Value = "Teleost Fish"
$ws.Range("A182").Value = "558d8758ae62abe36b1507ce2094ef7c"
$ws.Range("B182").Value = "Ammodytes dubius"
$ws.Range("C182").Value = "Northern sand lance"
$ws.Range("D182").Value = "Teleost Fish"
$ws.Range("A186").Value = "ee3c408644b66e62dde706ff463f359a"
$ws.Range("B186").Value = "Ammodytes dubius"
$ws.Range("C186").Value = "Northern sand lance"
$ws.Range("D186").Value = "Teleost Fish"
$ws.Range("A187").Value = "86c340102750abe5f2a75f3d5501b55d"
$ws.Range("B187").Value = "Menidia beryllina"
$ws.Range("C187").Value = "Inland silverside"
$ws.Range("D187").Value = "Teleost Fish"
$ws.Range("A190").Value = "0ad9142dc74ab0ef2021cfff48d4194d"
$ws.Range("B190").Value = "Brevoortia tyrannus"
$ws.Range("C190").Value = "Atlantic menhaden"
$ws.Range("D190").Value = "Teleost Fish"
$ws.Range("A191").Value = "8830d0cf4452e1cd0f9a6552b48b2b40"
$ws.Range("B191").Value = "Apeltes quadracus"
$ws.Range("C191").Value = "Fourspine stickleback"
$ws.Range("D191").Value = "Teleost Fish"
$ws.Range("A192").Value = "731abf4fa491ab03dd796729de5ab3eb"
$ws.Range("B192").Value = "Pseudopleuronectes americanus or Myzopsetta ferruginea"
$ws.Range("C192").Value = "Winter or Yellowtail flounder"
$ws.Range("D192").Value = "Teleost Fish"
$ws.Range("A195").Value = "f2e15a0b398b704a888c965d3b49035b"
$ws.Range("B195").Value = "Fundulus heteroclitus"
$ws.Range("C195").Value = "Mummichog"
$ws.Range("D195").Value = "Teleost Fish"
$ws.Range("A205").Value = "e468b57f39f048ada7562924022dc516"
$ws.Range("B205").Value = "Homo sapiens"
$ws.Range("C205").Value = "Human"
$ws.Range("D205").Value = "Human"
$ws.Range("A206").Value = "1533469db84e906a7d07208d202f0b61"
$ws.Range("B206").Value = "Fundulus heteroclitus"
$ws.Range("C206").Value = "Mummichog"
$ws.Range("D206").Value = "Teleost Fish"
$ws.Range("A216").Value = "5b2278535af7a77c15966bc43d0188bd"
$ws.Range("B216").Value = "Brevoortia tyrannus"
$ws.Range("C216").Value = "Atlantic menhaden"
$ws.Range("D216").Value = "Teleost Fish"
$ws.Range("A217").Value = "6a83eb23e34e01773abb7d038e38c583"
$ws.Range("B217").Value = "Pseudopleuronectes americanus or Myzopsetta ferruginea"
$ws.Range("C217").Value = "Winter or Yellowtail flounder"
$ws.Range("D217").Value = "Teleost Fish"
$ws.Range("A218").Value = "bdb87097756f45aa57e56f1d9f456f26"
$ws.Range("B218").Value = "Larus sp"
$ws.Range("C218").Value = "Great black backed gull and other gulls"
$ws.Range("D218").Value = "Bird"
$ws.Range("A220").Value = "cb17be39fabe38eb2368ba0635321393"
$ws.Range("B220").Value = "Tautoga onitis"
$ws.Range("C220").Value = "Tautog"
$ws.Range("D220").Value = "Teleost Fish"
$ws.Range("A221").Value = "c0b18824ab60460cd31eed51f737f882"
$ws.Range("B221").Value = "Homo sapiens"
$ws.Range("C221").Value = "Human"
$ws.Range("D221").Value = "Human"
$ws.Range("A225").Value = "879319f127f42872ba2daeb54fc4135a"
$ws.Range("B225").Value = "Gasterosteus aculeatus"
$ws.Range("C225").Value = "Threespined stickleback"
$ws.Range("D225").Value = "Teleost Fish"
$ws.Range("A226").Value = "f937641d91db232cb7180be9e04fb9e0"
$ws.Range("B226").Value = "Apeltes quadracus"
$ws.Range("C226").Value = "Fourspine stickleback"
$ws.Range("D226").Value = "Teleost Fish"
$ws.Range("A227").Value = "977b02be79d865979e54848db649eaf0"
$ws.Range("B227").Value = "Ammodytes dubius"
$ws.Range("C227").Value = "Northern sand lance"
$ws.Range("D227").Value = "Teleost Fish"
$ws.Range("A228").Value = "df263dae379496c7e522db8a7dbc01c9"
$ws.Range("B228").Value = "Scomber scombrus"
$ws.Range("C228").Value = "Atlantic mackerel"
$ws.Range("D228").Value = "Teleost Fish"
$ws.Range("A229").Value = "4450a6fa10b56881617cff33c5585aa8"
$ws.Range("B229").Value = "Fundulus heteroclitus"
$ws.Range("C229").Value = "Mummichog"
$ws.Range("D229").Value = "Teleost Fish"
$ws.Range("A230").Value = "f5ca5d430f1b145903b92fc335a4bafd"
$ws.Range("B230").Value = "Lucania parva"
$ws.Range("C230").Value = "Rainwater killifish"
$ws.Range("D230").Value = "Teleost Fish"
$ws.Range("A231").Value = "9c8a7b893d0fdaf8c1c89606cfce1c08"
$ws.Range("B231").Value = "Enchelyopus cimbrius"
$ws.Range("C231").Value = "Fourbeard rockling"
$ws.Range("D231").Value = "Teleost Fish"
$ws.Range("A232").Value = "c73c11c0d8b73d825099e12aaaefb637"
$ws.Range("B232").Value = "Fundulus heteroclitus"
$ws.Range("C232").Value = "Mummichog"
$ws.Range("D232").Value = "Teleost Fish"
$ws.Range("A233").Value = "c972de9c10572043855aaca4a4da68f4"
$ws.Range("B233").Value = "Ammodytes dubius"
$ws.Range("C233").Value = "Northern sand lance"
$ws.Range("D233").Value = "Teleost Fish"
$ws.Range("A234").Value = "243cc9259a8d104346a5dd517ca99499"
$ws.Range("B234").Value = "Ammodytes americanus"
$ws.Range("C234").Value = "American sand lance"
$ws.Range("D234").Value = "Teleost Fish"
$ws.Range("A244").Value = "f753730afbaa726c79bd991f32ea9778"
$ws.Range("B244").Value = "Brevoortia tyrannus"
$ws.Range("C244").Value = "Atlantic menhaden"
$ws.Range("D244").Value = "Teleost Fish"
$ws.Range("A245").Value = "5432a6e652c21bb79c110c1179832080"
$ws.Range("B245").Value = "Clupeidae sp"
$ws.Range("C245").Value = "Atlantic menhaden or River herrings"
$ws.Range("D245").Value = "Teleost Fish"
$ws.Range("A246").Value = "279fde05b5aed4bfe15ab39776ff82ba"
$ws.Range("B246").Value = "Brevoortia tyrannus"
$ws.Range("C246").Value = "Atlantic menhaden"
$ws.Range("D246").Value = "Teleost Fish"
$ws.Range("A248").Value = "4db280926cca07cc86b0e098513d9cc0"
$ws.Range("B248").Value = "Ammodytes dubius"
$ws.Range("C248").Value = "Northern sand lance"
$ws.Range("D248").Value = "Teleost Fish"
$ws.Range("A249").Value = "88065f0fd14ae3b76fc1a87f8df6ef2d"
$ws.Range("B249").Value = "Pseudopleuronectes americanus or Myzopsetta ferruginea"
$ws.Range("C249").Value = "Winter or Yellowtail flounder"
$ws.Range("D249").Value = "Teleost Fish"
$ws.Range("A251").Value = "fc6d040e1564a91e1c6d67e1e32b9022"
$ws.Range("B251").Value = "Rattus norvegicus"
$ws.Range("C251").Value = "Norway rat"
$ws.Range("D251").Value = "Mammal"
$ws.Range("A258").Value = "14bd3bb11b9a6c641ad60556bf6141d0"
$ws.Range("B258").Value = "Pseudopleuronectes americanus or Myzopsetta ferruginea"
$ws.Range("C258").Value = "Winter or Yellowtail flounder"
$ws.Range("D258").Value = "Teleost Fish"
$ws.Range("A259").Value = "8103469b2716037f1cc4ce8959ae0081"
$ws.Range("B259").Value = "Menidia menidia"
$ws.Range("C259").Value = "Atlantic silverside"
$ws.Range("D259").Value = "Teleost Fish"
$ws.Range("A262").Value = "f4d5447013c09b659b99d47459de2042"
$ws.Range("B262").Value = "Fundulus heteroclitus"
$ws.Range("C262").Value = "Mummichog"
$ws.Range("D262").Value = "Teleost Fish"
$ws.Range("A263").Value = "29d8e064f48ae7211c9fba32872b36f9"
$ws.Range("B263").Value = "Apeltes quadracus"
$ws.Range("C263").Value = "Fourspine stickleback"
$ws.Range("D263").Value = "Teleost Fish"
$ws.Range("A264").Value = "53cfac0a209f1dbdaf758a75c84df7d6"
$ws.Range("B264").Value = "Homo sapiens"
$ws.Range("C264").Value = "Human"
$ws.Range("D264").Value = "Human"
$ws.Range("A265").Value = "9dc70d19a67c006232234c9bcbbab33f"
$ws.Range("B265").Value = "Clangula hyemalis or other Anatidae sp"
$ws.Range("C265").Value = "Long tailed duck or other ducks"
$ws.Range("D265").Value = "Bird"
$ws.Range("A267").Value = "daa05108bed6292fbd2eedef6214fdff"
$ws.Range("B267").Value = "Fundulus heteroclitus"
$ws.Range("C267").Value = "Mummichog"
$ws.Range("D267").Value = "Teleost Fish"
$ws.Range("A269").Value = "4c451c35f94e737edb8816211106c35d"
$ws.Range("B269").Value = "Fundulus heteroclitus"
$ws.Range("C269").Value = "Mummichog"
$ws.Range("D269").Value = "Teleost Fish"
$ws.Range("A270").Value = "b61cc19540f9627af5070110979ebf91"
$ws.Range("B270").Value = "Ammodytes dubius"
$ws.Range("C270").Value = "Northern sand lance"
$ws.Range("D270").Value = "Teleost Fish"
$ws.Range("A271").Value = "26a933bf07de2306f33cc95ec94e4b2f"
$ws.Range("B271").Value = "Pomoxis nigromaculatus"
$ws.Range("C271").Value = "Black crappie"
$ws.Range("D271").Value = "Teleost Fish"
$ws.Range("A273").Value = "93b36a6e82074114f7c4d90b6172dba2"
$ws.Range("B273").Value = "Fundulus heteroclitus"
$ws.Range("C273").Value = "Mummichog"
$ws.Range("D273").Value = "Teleost Fish"
$ws.Range("A274").Value = "1be2860881962b4dd3a0a7c6db14ca80"
$ws.Range("B274").Value = "Pseudopleuronectes americanus or Myzopsetta ferruginea"
$ws.Range("C274").Value = "Winter or Yellowtail flounder"
$ws.Range("D274").Value = "Teleost Fish"
$ws.Range("A278").Value = "b7b35bf53a25eef31602b3f785c925e9"
$ws.Range("B278").Value = "Anguilla rostrata"
$ws.Range("C278").Value = "American eel"
$ws.Range("D278").Value = "Teleost Fish"
$ws.Range("A280").Value = "5cc8cf140b434a5aba3a1dad41339918"
$ws.Range("B280").Value = "Anguilla rostrata"
$ws.Range("C280").Value = "American eel"
$ws.Range("D280").Value = "Teleost Fish"
$ws.Range("A283").Value = "29df87a23b45339e281dc7b390d16860"
$ws.Range("B283").Value = "Apeltes quadracus"
$ws.Range("C283").Value = "Fourspine stickleback"
$ws.Range("D283").Value = "Teleost Fish"
$ws.Range("A284").Value = "8dde047966dfd43b699a5ca7122e55d7"
$ws.Range("B284").Value = "Pseudopleuronectes americanus or Myzopsetta ferruginea"
$ws.Range("C284").Value = "Winter or Yellowtail flounder"
$ws.Range("D284").Value = "Teleost Fish"
$ws.Range("A285").Value = "118de0da9053ad27ad0e3c1e136454d9"
$ws.Range("B285").Value = "Ammodytes dubius"
$ws.Range("C285").Value = "Northern sand lance"
$ws.Range("D285").Value = "Teleost Fish"
$ws.Range("A287").Value = "ed57094384d61f9a2dbc0c0e4ff6fb8a"
$ws.Range("B287").Value = "Lepomis gibbosus"
$ws.Range("C287").Value = "Pumpkinseed"
$ws.Range("D287").Value = "Teleost Fish"
$ws.Range("A288").Value = "d8de668ed4c19b2b2ea845f3db18ae79"
$ws.Range("B288").Value = "Cottidae sp"
$ws.Range("C288").Value = "Sculpins"
$ws.Range("D288").Value = "Teleost Fish"
$ws.Range("A289").Value = "7bb99841c9c9b86e238cbcddc0e16567"
$ws.Range("B289").Value = "Micropterus salmoides"
$ws.Range("C289").Value = "Largemouth bass"
$ws.Range("D289").Value = "Teleost Fish"
$ws.Range("A292").Value = "1dd0f1ca2adf649d8cba813ea6e43de2"
$ws.Range("B292").Value = "Paralichthys dentatus"
$ws.Range("C292").Value = "Summer flounder"
$ws.Range("D292").Value = "Teleost Fish"
$ws.Range("A293").Value = "50bddde558bebcd7fa8dbf6542ea44fe"
$ws.Range("B293").Value = "Fundulus heteroclitus"
$ws.Range("C293").Value = "Mummichog"
$ws.Range("D293").Value = "Teleost Fish"
$ws.Range("A295").Value = "92693323f831e69117617606814ae81f"
$ws.Range("B295").Value = "Coryphaena hippurus"
$ws.Range("C295").Value = "Mahi mahi"
$ws.Range("D295").Value = "Teleost Fish"
$ws.Range("A297").Value = "0d6e610cd1019f50d693803e46db364f"
$ws.Range("B297").Value = "Pseudopleuronectes americanus or Myzopsetta ferruginea"
$ws.Range("C297").Value = "Winter or Yellowtail flounder"
$ws.Range("D297").Value = "Teleost Fish"
$ws.Range("A298").Value = "bbd6b723329db44753870a56d15bdbd6"
$ws.Range("B298").Value = "Sebastes fasciatus"
$ws.Range("C298").Value = "Acadian redfish"
$ws.Range("D298").Value = "Teleost Fish"
$ws.Range("A303").Value = "9f4285ab8775db6b862ee4fb416f0f5d"
$ws.Range("B303").Value = "Melospiza melodia or Spizella passerina"
$ws.Range("C303").Value = "Song sparrow or Chipping sparrow"
$ws.Range("D303").Value = "Bird"
$ws.Range("A304").Value = "0a6108b3c6bbca90164970efbea23261"
$ws.Range("B304").Value = "Apeltes quadracus"
$ws.Range("C304").Value = "Fourspine stickleback"
$ws.Range("D304").Value = "Teleost Fish"
$ws.Range("A305").Value = "901fc1f68af659cc3f6678c6a7396845"
$ws.Range("B305").Value = "Brevoortia tyrannus"
$ws.Range("C305").Value = "Atlantic menhaden"
$ws.Range("D305").Value = "Teleost Fish"
$ws.Range("A307").Value = "5d4b77f374dcda6b5f48e88cc2b9664b"
$ws.Range("B307").Value = "Unassigned"
$ws.Range("C307").Value = "Unassigned"
$ws.Range("D307").Value = "Unassigned"
$ws.Range("A308").Value = "7950b1078efc076defba9c936b970ef7"
$ws.Range("B308").Value = "Homo sapiens"
$ws.Range("C308").Value = "Human"
$ws.Range("D308").Value = "Human"
$ws.Range("A309").Value = "8ea2a9236bef33ba65acfc82e6947942"
$ws.Range("B309").Value = "Sciurus carolinensis"
$ws.Range("C309").Value = "Gray squirrel"
$ws.Range("D309").Value = "Mammal"
$ws.Range("A310").Value = "c1f17b3dc22ac71ee83288f654c93bb3"
$ws.Range("B310").Value = "Lontra canadensis"
$ws.Range("C310").Value = "River otter"
$ws.Range("D310").Value = "Mammal"
$ws.Range("A311").Value = "a4e8997c6347c55b72f81e0accce0c37"
$ws.Range("B311").Value = "Pseudopleuronectes americanus or Myzopsetta ferruginea"
$ws.Range("C311").Value = "Winter or Yellowtail flounder"
$ws.Range("D311").Value = "Teleost Fish"
$ws.Range("A313").Value = "7eef5797ad87b51600785f22606c70de"
$ws.Range("B313").Value = "Brevoortia tyrannus"
$ws.Range("C313").Value = "Atlantic menhaden"
$ws.Range("D313").Value = "Teleost Fish"
$ws.Range("A314").Value = "033531a8711295f5cf38c1111629eb77"
$ws.Range("B314").Value = "Anguilla rostrata"
$ws.Range("C314").Value = "American eel"
$ws.Range("D314").Value = "Teleost Fish"
$ws.Range("A316").Value = "16d55edf1062cb60bf8a36a1da3212b5"
$ws.Range("B316").Value = "Esox americanus or niger"
$ws.Range("C316").Value = "Grass or chain pickerel"
$ws.Range("D316").Value = "Teleost Fish"
$ws.Range("A317").Value = "5b1dbdcc719bcfd9ea209ec7d9ecd075"
$ws.Range("B317").Value = "Pseudopleuronectes americanus or Myzopsetta ferruginea"
$ws.Range("C317").Value = "Winter or Yellowtail flounder"
$ws.Range("D317").Value = "Teleost Fish"
$ws.Range("A319").Value = "d1de955bd9480b3f0c70f78caec5a443"
$ws.Range("B319").Value = "Homo sapiens"
$ws.Range("C319").Value = "Human"
$ws.Range("D319").Value = "Human"
$ws.Range("A320").Value = "bfab25a003878187c8038ee55fdb7a53"
$ws.Range("B320").Value = "Gavia immer"
$ws.Range("C320").Value = "Common loon"
$ws.Range("D320").Value = "Bird"
$ws.Range("A321").Value = "fb3bb0a4483dcfbc39e8b7ccf8196749"
$ws.Range("B321").Value = "Homo sapiens"
$ws.Range("C321").Value = "Human"
$ws.Range("D321").Value = "Human"
$ws.Range("A322").Value = "9ed3306f1d1dfb81749820128e325abc"
$ws.Range("B322").Value = "Clupeidae sp"
$ws.Range("C322").Value = "Atlantic menhaden or River herrings"
$ws.Range("D322").Value = "Teleost Fish"
$ws.Range("A327").Value = "9f2355fd161fec0177a83045e771a239"
$ws.Range("B327").Value = "Decapterus punctatus"
$ws.Range("C327").Value = "Round scad"
$ws.Range("D327").Value = "Teleost Fish"
$ws.Range("A328").Value = "032747b5e01bbdc74f3e3e59d9c2275f"
$ws.Range("B328").Value = "Unassigned"
$ws.Range("C328").Value = "Unassigned"
$ws.Range("D328").Value = "Unassigned"
$ws.Range("A336").Value = "856a99622e9c49fc86cdcf3ddbefcfd1"
$ws.Range("B336").Value = "Gobiosoma ginsburgi"
$ws.Range("C336").Value = "Seaboard goby"
$ws.Range("D336").Value = "Teleost Fish"
$ws.Range("A337").Value = "38e4381a618398d035b19d8c47e8bfa3"
$ws.Range("B337").Value = "Ammodytes dubius"
$ws.Range("C337").Value = "Northern sand lance"
$ws.Range("D337").Value = "Teleost Fish"
$ws.Range("A338").Value = "a1f66fcba0bd12eed7a60901bd8c6010"
$ws.Range("B338").Value = "Pseudopleuronectes americanus or Myzopsetta ferruginea"
$ws.Range("C338").Value = "Winter or Yellowtail flounder"
$ws.Range("D338").Value = "Teleost Fish"
$ws.Range("A340").Value = "5cde257b7febb75f7c9848a21bfe18cd"
$ws.Range("B340").Value = "Fundulus heteroclitus"
$ws.Range("C340").Value = "Mummichog"
$ws.Range("D340").Value = "Teleost Fish"
$ws.Range("A341").Value = "d3b57d4cf93def7c41d4b9baced940d9"
$ws.Range("B341").Value = "Canis lupus"
$ws.Range("C341").Value = "Dog"
$ws.Range("D341").Value = "Livestock"
$ws.Range("A343").Value = "60fb36f888cbf7e4639c1bb98f0adc57"
$ws.Range("B343").Value = "Menidia menidia"
$ws.Range("C343").Value = "Atlantic silverside"
$ws.Range("D343").Value = "Teleost Fish"
$ws.Range("A345").Value = "e7f90ff8c7b97da66ce6d940d857e8e4"
$ws.Range("B345").Value = "Trachurus lathami"
$ws.Range("C345").Value = "Rough scad"
$ws.Range("D345").Value = "Teleost Fish"
$ws.Range("A346").Value = "bcf13bc540e00c02358754d8a1b40a9c"
$ws.Range("B346").Value = "Fundulus heteroclitus"
$ws.Range("C346").Value = "Mummichog"
$ws.Range("D346").Value = "Teleost Fish"
$ws.Range("A347").Value = "7fef2f8e6a8bee56528216dfc05f0d81"
$ws.Range("B347").Value = "Homo sapiens"
$ws.Range("C347").Value = "Human"
$ws.Range("D347").Value = "Human"
$ws.Range("A353").Value = "1a9a786e8451eec71300762a5398f4d0"
$ws.Range("B353").Value = "Anguilla rostrata"
$ws.Range("C353").Value = "American eel"
$ws.Range("D353").Value = "Teleost Fish"
$ws.Range("A354").Value = "0e3aec812235602fac414c57ef969f1e"
$ws.Range("B354").Value = "Peromyscus leucopus"
$ws.Range("C354").Value = "Deer mouse"
$ws.Range("D354").Value = "Mammal"
$ws.Range("A355").Value = "35a8484c22fbf1df676003af6ec52a29"
$ws.Range("B355").Value = "Pseudopleuronectes americanus or Myzopsetta ferruginea"
$ws.Range("C355").Value = "Winter or Yellowtail flounder"
$ws.Range("D355").Value = "Teleost Fish"
$ws.Range("A356").Value = "841570d820eaab46bb7c7b3b7db8fba3"
$ws.Range("B356").Value = "Ctenogobius boleosoma"
$ws.Range("C356").Value = "Darter goby"
$ws.Range("D356").Value = "Teleost Fish"
$ws.Range("A357").Value = "b45f5a79f080475521114b63aa9bd7d1"
$ws.Range("B357").Value = "Thunnus sp"
$ws.Range("C357").Value = "Tuna sp"
$ws.Range("D357").Value = "Teleost Fish"
$ws.Range("A358").Value = "5d6bb2fed75f92f01e645ffc80d17d36"
$ws.Range("B358").Value = "Anchoa mitchilli"
$ws.Range("C358").Value = "Bay anchovy"
$ws.Range("D358").Value = "Teleost Fish"
$ws.Range("A359").Value = "cc23248fc9f0058810041c6090c99461"
$ws.Range("B359").Value = "Unassigned"
$ws.Range("C359").Value = "Unassigned"
$ws.Range("D359").Value = "Unassigned"
$ws.Range("A364").Value = "02af37069efe14191838e3ded56eba42"
$ws.Range("B364").Value = "Unassigned"
$ws.Range("C364").Value = "Unassigned"
$ws.Range("D364").Value = "Unassigned"
$ws.Range("A365").Value = "0dc4976a75f5949215d7cf1f1a2994f6"
$ws.Range("B365").Value = "Unassigned"
$ws.Range("C365").Value = "Unassigned"
$ws.Range("D365").Value = "Unassigned"
$ws.Range("A366").Value = "dedc23ac2bedc9a49f0d9d5f2092f08a"
$ws.Range("B366").Value = "Sternotherus carinatus"
$ws.Range("C366").Value = "Razor-backed musk turtle"
$ws.Range("D366").Value = "Reptile"
$ws.Range("A368").Value = "072eff78dfd4d6ee7d89b57c63c82827"
$ws.Range("B368").Value = "Cottidae sp"
$ws.Range("C368").Value = "Sculpins"
$ws.Range("D368").Value = "Teleost Fish"
$ws.Range("A369").Value = "f3ca4711d94796503be4e5bbb8a6705b"
$ws.Range("B369").Value = "Homo sapiens"
$ws.Range("C369").Value = "Human"
$ws.Range("D369").Value = "Human"
$ws.Range("A376").Value = "66ce172600ad7ff5f35d5dfc0bab87d3"
$ws.Range("B376").Value = "Homo sapiens"
$ws.Range("C376").Value = "Human"
$ws.Range("D376").Value = "Human"
$ws.Range("A377").Value = "9b1cbc97eff7dc6af2c100d85526140f"
$ws.Range("B377").Value = "Homo sapiens"
$ws.Range("C377").Value = "Human"
$ws.Range("D377").Value = "Human"
$ws.Range("A380").Value = "798cab8e6a1a556a317f78cfa6bab8ac"
$ws.Range("B380").Value = "Felis catus"
$ws.Range("C380").Value = "Cat"
$ws.Range("D380").Value = "Livestock"
$ws.Range("A381").Value = "97444d2388851f1d71afeb95125b4898"
$ws.Range("B381").Value = "Unassigned"
$ws.Range("C381").Value = "Unassigned"
$ws.Range("D381").Value = "Unassigned"
$ws.Range("A387").Value = "41714252fc55bbced79f657eb2b8805a"
$ws.Range("B387").Value = "Passer domesticus"
$ws.Range("C387").Value = "House sparrow"
$ws.Range("D387").Value = "Bird"
$ws.Range("A389").Value = "3b578403acdc73dd077d282c96f9541f"
$ws.Range("B389").Value = "Ophidion marginatum"
$ws.Range("C389").Value = "Striped cusk-eel"
$ws.Range("D389").Value = "Teleost Fish"
$ws.Range("A393").Value = "d1af29b8548fccd9f3eada2b18f0eac9"
$ws.Range("B393").Value = "Unassigned"
$ws.Range("C393").Value = "Unassigned"
$ws.Range("D393").Value = "Unassigned"
$ws.Range("A394").Value = "803a43fa7cb74bb51f36ab2949523bf2"
$ws.Range("B394").Value = "Homo sapiens"
$ws.Range("C394").Value = "Human"
$ws.Range("D394").Value = "Human"
$ws.Range("A397").Value = "812ed0386e2a4869a21da5634665548d"
$ws.Range("B397").Value = "Caranx hippos"
$ws.Range("C397").Value = "Crevalle jack"
$ws.Range("D397").Value = "Teleost Fish"
$ws.Range("A404").Value = "a6bf1361741b5eda21b4d05f18f04a90"
$ws.Range("B404").Value = "Unassigned"
$ws.Range("C404").Value = "Unassigned"
$ws.Range("D404").Value = "Unassigned"
$ws.Range("A405").Value = "2d0506c060ee125f6608b52f22e598b1"
$ws.Range("B405").Value = "Homo sapiens"
$ws.Range("C405").Value = "Human"
$ws.Range("D405").Value = "Human"
$ws.Range("A407").Value = "9e218ddde05826daea9943de26124674"
$ws.Range("B407").Value = "Sturnus vulgaris"
$ws.Range("C407").Value = "Common starling"
$ws.Range("D407").Value = "Bird"
$ws.Range("A408").Value = "f6b3a673a06591a2d5a8936584e64754"
$ws.Range("B408").Value = "Tautogolabrus adspersus"
$ws.Range("C408").Value = "Cunner"
$ws.Range("D408").Value = "Teleost Fish"
$ws.Range("A409").Value = "75fedd0fa34e3ac2514601e68b613736"
$ws.Range("B409").Value = "Strongylura marina"
$ws.Range("C409").Value = "Atlantic needlefish"
$ws.Range("D409").Value = "Teleost Fish"
$ws.Range("A411").Value = "037bd6992d173dfbcd22d76af622fa5b"
$ws.Range("B411").Value = "Malaclemys terrapin"
$ws.Range("C411").Value = "Diamondback terrapin"
$ws.Range("D411").Value = "Reptile"
$ws.Range("A412").Value = "08bd987bd944513cc896ab3b3c3eed38"
$ws.Range("B412").Value = "Hippoglossina oblonga"
$ws.Range("C412").Value = "Fourspot flounder"
$ws.Range("D412").Value = "Teleost Fish"
$ws.Range("A413").Value = "0437f2363acf453d6291ceda4abba683"
$ws.Range("B413").Value = "Ameiurus nebulosus"
$ws.Range("C413").Value = "Brown bullhead"
$ws.Range("D413").Value = "Teleost Fish"
$ws.Range("A414").Value = "956d5064b9d6c222e19d75e231925e18"
$ws.Range("B414").Value = "Homo sapiens"
$ws.Range("C414").Value = "Human"
$ws.Range("D414").Value = "Human"
$ws.Range("A415").Value = "d964add43fe0c3212cbe19a066dc2a13"
$ws.Range("B415").Value = "Unassigned"
$ws.Range("C415").Value = "Unassigned"
$ws.Range("D415").Value = "Unassigned"
$ws.Range("A416").Value = "7f400300a06f165c23af04aa4e4c790c"
$ws.Range("B416").Value = "Enchelyopus cimbrius"
$ws.Range("C416").Value = "Fourbeard rockling"
$ws.Range("D416").Value = "Teleost Fish"
$ws.Range("A417").Value = "ced5f183dc83bf9a3831984cefa3b3a1"
$ws.Range("B417").Value = "Homo sapiens"
$ws.Range("C417").Value = "Human"
$ws.Range("D417").Value = "Human"
$ws.Range("A418").Value = "29ae99676d91135240ab43f0184c5909"
$ws.Range("B418").Value = "Unassigned"
$ws.Range("C418").Value = "Unassigned"
$ws.Range("D418").Value = "Unassigned"
$ws.Range("A419").Value = "df6846357baa6a9fecd66b4a1ba513a8"
$ws.Range("B419").Value = "Canis lupus"
$ws.Range("C419").Value = "Dog"
$ws.Range("D419").Value = "Livestock"
$ws.Range("A420").Value = "39c09623e77e6cb1f69a264089e6256c"
$ws.Range("B420").Value = "Blarina brevicauda"
$ws.Range("C420").Value = "Northern short tailed shrew"
$ws.Range("D420").Value = "Mammal"
